$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3) appended below existing data
$ws.Range("A3").Value = 42632.88212962963
$ws.Range("A3").NumberFormat = "m/d/yyyy h:mm:ss"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "Neutral"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 8312
$ws.Range("F3").Value = 1284
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 38
$ws.Range("I3").Value = 75
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 10254
$ws.Range("L3").Value = 151
$ws.Range("M3").Value = 95
$ws.Range("N3").Value = 15
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = "Noun"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1.77
$ws.Range("S3").Value = 0.1132
$ws.Range("S3").NumberFormat = "0.00%"
$ws.Range("T3").Value = -4.05
$ws.Range("U3").Value = 5.85
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0
